$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update the two price cells
$ws.Range("D29").Value = 23711
$ws.Range("D30").Value = 26200
